# Fruta / hortaliza, semanal
# A new weekly record is prepended to the "Pomelo" price list: insert a new
# row at position 106 (shifting all subsequent rows down by one, 106-189 -> 107-190)
# and populate the new row 106 with this week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 106; existing rows 106..189 shift down to 107..190
$ws.Rows(106).Insert()

# Populate the newly inserted row 106 with the new weekly record
$ws.Range("A106").Value = 10
$ws.Range("B106").Value = "Vega Modelo de Temuco"
$ws.Range("C106").Value = "La Araucanía"
$ws.Range("D106").Value = 44589
$ws.Range("E106").Value = 9
$ws.Range("F106").Value = "Fruta"
$ws.Range("G106").Value = 100102
$ws.Range("H106").Value = "Cítricos"
$ws.Range("I106").Value = 100102006
$ws.Range("J106").Value = "Pomelo"
$ws.Range("K106").Value = "Start Ruby"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 50
$ws.Range("N106").Value = 14000
$ws.Range("O106").Value = 14000
$ws.Range("P106").Value = 14000
$ws.Range("Q106").Value = "$/bandeja 15 kilos granel"
$ws.Range("R106").Value = "Región de O'Higgins"
$ws.Range("S106").Value = 933
$ws.Range("T106").Value = 15
